$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.02,1.032330411378405,1.050090301998843,1.042570697001373,1.055438703108253,1.0371631416238,1.037461224153471,1.052845545953301,1.045347010076581,1.058179159405954,1.038934538090235)
    3 = @(1.02,1.033427835524426,1.050762338969189,1.043506117794467,1.05636636981964,1.037278048912751,1.038200146232684,1.053330321739099,1.046092963146672,1.058919976155134,1.039674509523555)
    4 = @(1.02,1.034138530754517,1.051196767554782,1.0441120952042,1.056966943621954,1.037350882755665,1.038678323185055,1.053642899056258,1.046575716871396,1.059399016363987,1.040153365542046)
    5 = @(1.02,1.034437447968713,1.051379298188842,1.044367014342205,1.057219498245462,1.03738113788672,1.038879359162553,1.053774040804262,1.0467786832126,1.059600327790475,1.040354687013692)
    6 = @(1.02,1.034487645751206,1.051409939791293,1.04440982611282,1.057261907548494,1.037386196481075,1.038913114585986,1.053796044461502,1.046812763114491,1.059634124315749,1.040388490373699)
    7 = @(1.02,1.034142524347069,1.051199206945328,1.044115500793982,1.056970317982256,1.037351288458086,1.038681009397284,1.053644652424491,1.04657842885283,1.059401706601255,1.040156055569006)
    8 = @(1.02,1.032701169209443,1.050317506944462,1.042886681928638,1.055752147243461,1.03720228920964,1.037710936794531,1.05300960627827,1.045599093090227,1.058429587087024,1.039184605351887)
    9 = @(1.02,1.030165820607669,1.048760658606318,1.040726722877036,1.053608010547559,1.036928127224596,1.036001905653516,1.051882163761616,1.043873958183422,1.05671419870409,1.037473147190624)
    10 = @(1.02,1.028478608655303,1.047720713811343,1.039290409332847,1.052180290707649,1.036737584309498,1.034862814517298,1.051124951774597,1.042724292709452,1.055569057923208,1.036332438414336)
    11 = @(1.02,1.027748741182887,1.047269940173344,1.03866934637367,1.051562490198061,1.036653241866687,1.034369641595869,1.050795760910345,1.042226582344272,1.055072843840579,1.035838565130791)
    12 = @(1.02,1.027477741726748,1.047102433571698,1.038438787332832,1.051333074460989,1.036621638027369,1.034186464618654,1.050673288627707,1.04204172642002,1.054888474558249,1.035655128021258)
    13 = @(1.02,1.027535867246388,1.047138367413667,1.038488237073478,1.051382282053213,1.036628429607254,1.034225756272772,1.050699568205542,1.042081377917501,1.054928024791779,1.035694475474032)
    14 = @(1.02,1.027726338126387,1.047256095428369,1.038650285607281,1.051543525332924,1.036650635096328,1.034354499936599,1.050785641307864,1.042211301769434,1.055057604899884,1.035823401968628)
    15 = @(1.02,1.027843707534863,1.04732862241897,1.038750146487654,1.051642881002556,1.036664280167525,1.034433824410083,1.050838647830657,1.042291354259579,1.055137436448002,1.035902839091964)
    16 = @(1.02,1.028527062412572,1.047750620385672,1.039331645626914,1.052221300853678,1.036743143202302,1.034895546085949,1.051146771520172,1.042757326305508,1.055601982497138,1.036365216465569)
    17 = @(1.02,1.028955901783295,1.048015203616552,1.039696638128599,1.052584239277762,1.036792120714354,1.035185188404038,1.051339698625718,1.043049646009573,1.0558932840573,1.036655270108972)
    18 = @(1.02,1.029206104714445,1.048169485023875,1.039909616006106,1.052795974763996,1.036820511274134,1.03535413774333,1.051452102973382,1.043220160992027,1.056063160479481,1.036824459375705)
    19 = @(1.02,1.029291429022564,1.048222083230584,1.039982250174278,1.052868177785676,1.036830161649264,1.035411746091581,1.051490408414312,1.043278303849623,1.0561210780112,1.036882149534419)
    20 = @(1.02,1.028909884327116,1.047986821036713,1.039657469201787,1.05254529531153,1.03678688421486,1.035154111914979,1.051319012477187,1.04301828184995,1.055862033754869,1.036624149487735)
    21 = @(1.02,1.027670246280913,1.047221429360523,1.0386025627207,1.051496041461786,1.036644103728464,1.034316587874979,1.050760300329817,1.042173041991422,1.055019448251314,1.035785436067532)
    22 = @(1.02,1.026891448734361,1.046739797599631,1.037940061677968,1.050836698239248,1.0365527394376,1.033790057796625,1.050407881590884,1.041641698181632,1.054489373493873,1.035258158256073)
    23 = @(1.02,1.027304246132738,1.046995157086241,1.038291193710585,1.051186193676271,1.036601324146169,1.03406917610713,1.050594812566451,1.041923364702069,1.054770405030959,1.035537672946594)
    24 = @(1.02,1.028930677423901,1.04799964604359,1.039675167694474,1.05256289229109,1.036789250911119,1.035168154024489,1.051328360049055,1.04303245393318,1.055876154529084,1.03663821153865)
    25 = @(1.02,1.030820736103857,1.049163508562618,1.041284482346106,1.054162025722583,1.037000376446027,1.036443686323376,1.052174624001602,1.044319874630843,1.057157944906388,1.037915555239708)
}

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

Write-Output "Updated vm_pu values for case with 380 kV"